# Updates the coinranking.com crypto price/volume snapshot table on Sheet1.
# Mirrors an automated "refresh data" run: most cells just get a literal new
# text value, but two pairs of rows (39/40 and 47/48) were re-ranked and swap
# their Coin/Link/Price/Volume content while the leading rank index (col A)
# stays put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. "29.338.06" or "0.9989").
# Excel auto-detects genuinely numeric-looking text and would silently convert
# it to a Number cell, which would lose the original text formatting. For those
# values we momentarily force Text format, write the value, then drop back to
# the worksheet default style so no stray per-cell formatting is left behind.
function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "29.338.06"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3
$ws.Range("D3").Value = "1.840.07"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4
Set-TextValue $ws "D4" "0.9989"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
Set-TextValue $ws "D5" "239.09"

# Row 6
Set-TextValue $ws "D6" "0.6272"
$ws.Range("E6").Value = "  -0.88%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
Set-TextValue $ws "D8" "0.07432"
$ws.Range("E8").Value = "  -0.89%  "

# Row 9
Set-TextValue $ws "D9" "0.2887"
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
Set-TextValue $ws "D10" "24.89"
$ws.Range("E10").Value = "  +1.92%  "

# Row 11
Set-TextValue $ws "D11" "0.07728"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("D12").Value = "1.838.39"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13
Set-TextValue $ws "D13" "4.956"
$ws.Range("E13").Value = "  -0.98%  "

# Row 14
Set-TextValue $ws "D14" "0.6742"
$ws.Range("E14").Value = "  -0.72%  "

# Row 15
Set-TextValue $ws "D15" "0.00001024"
$ws.Range("E15").Value = "  -0.59%  "

# Row 16
Set-TextValue $ws "D16" "81.49"
$ws.Range("E16").Value = "  -0.81%  "

# Row 17
Set-TextValue $ws "D17" "6.220"
$ws.Range("E17").Value = "  +1.16%  "

# Row 18
$ws.Range("D18").Value = "29.376.52"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
Set-TextValue $ws "D19" "228.60"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
Set-TextValue $ws "D20" "12.28"
$ws.Range("E20").Value = "  -0.41%  "

# Row 21
Set-TextValue $ws "D21" "0.9998"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
Set-TextValue $ws "D22" "7.327"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23
Set-TextValue $ws "D23" "1.001"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
Set-TextValue $ws "D24" "157.82"
$ws.Range("E24").Value = "  -0.76%  "

# Row 25
Set-TextValue $ws "D25" "8.458"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
Set-TextValue $ws "D26" "0.1341"
$ws.Range("E26").Value = "  -2.78%  "

# Row 27
Set-TextValue $ws "D27" "17.36"
$ws.Range("E27").Value = "  -1.29%  "

# Row 28
Set-TextValue $ws "D28" "0.07190"
$ws.Range("E28").Value = "  +12.97%  "

# Row 29
Set-TextValue $ws "D29" "1.459"
$ws.Range("E29").Value = "  +5.43%  "

# Row 30
Set-TextValue $ws "D30" "1.478"
$ws.Range("E30").Value = "  +0.14%  "

# Row 31
Set-TextValue $ws "D31" "4.035"
$ws.Range("E31").Value = "  -1.43%  "

# Row 32
Set-TextValue $ws "D32" "4.034"
$ws.Range("E32").Value = "  -0.36%  "

# Row 33
Set-TextValue $ws "D33" "1.819"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34
Set-TextValue $ws "D34" "1.137"
$ws.Range("E34").Value = "  -0.37%  "

# Row 35
Set-TextValue $ws "D35" "0.6940"
$ws.Range("E35").Value = "  -0.79%  "

# Row 36
Set-TextValue $ws "D36" "2.579"

# Row 37
Set-TextValue $ws "D37" "0.01839"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38
Set-TextValue $ws "D38" "2.804"
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D39" "6.813"
$ws.Range("E39").Value = "  +3.52%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.232.94"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
Set-TextValue $ws "D41" "0.9272"
$ws.Range("E41").Value = "  +2.25%  "

# Row 42
Set-TextValue $ws "D42" "0.9998"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43
$ws.Range("D43").Value = "2.011.02"
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
Set-TextValue $ws "D44" "100.35"
$ws.Range("E44").Value = "  -1.03%  "

# Row 45
Set-TextValue $ws "D45" "65.24"
$ws.Range("E45").Value = "  -1.69%  "

# Row 46
Set-TextValue $ws "D46" "0.00000000120"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D47" "1.705"
$ws.Range("E47").Value = "  +0.50%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D48" "6.938"
$ws.Range("E48").Value = "  -1.66%  "

# Row 49
Set-TextValue $ws "D49" "8.914"
$ws.Range("E49").Value = "  -1.33%  "

# Row 50
Set-TextValue $ws "D50" "0.1135"
$ws.Range("E50").Value = "  -3.65%  "

# Row 51
Set-TextValue $ws "D51" "0.3896"
$ws.Range("E51").Value = "  -1.09%  "
